$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.416.52'
$ws.Range("E2").Value = '  +2.56%  '

# Row 3
$ws.Range("D3").Value = '2.943.18'
$ws.Range("E3").Value = '  +2.41%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.97'
$ws.Range("E5").Value = '  +0.08%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.93'
$ws.Range("E6").Value = '  +4.94%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '2.943.10'
$ws.Range("E8").Value = '  +2.53%  '

# Row 9
$ws.Range("E9").Value = '  +3.36%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.98'
$ws.Range("E10").Value = '  +3.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.149'
$ws.Range("E11").Value = '  +9.42%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.434'
$ws.Range("E12").Value = '  +1.43%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  +7.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.20'
$ws.Range("E14").Value = '  -0.01%  '

# Row 15
$ws.Range("E15").Value = '  -0.95%  '

# Row 16
$ws.Range("D16").Value = '3.435.93'
$ws.Range("E16").Value = '  +2.57%  '

# Row 17
$ws.Range("D17").Value = '62.474.41'
$ws.Range("E17").Value = '  +2.72%  '

# Row 18
$ws.Range("D18").Value = '2.946.82'
$ws.Range("E18").Value = '  +2.49%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.64'
$ws.Range("E19").Value = '  +2.29%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '434.28'
$ws.Range("E20").Value = '  +2.56%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.42'
$ws.Range("E21").Value = '  +1.54%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.661'
$ws.Range("E22").Value = '  +1.43%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.95'
$ws.Range("E23").Value = '  +0.81%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.04'
$ws.Range("E24").Value = '  +6.14%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.05'
$ws.Range("E25").Value = '  +0.77%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.91'
$ws.Range("E26").Value = '  +5.41%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.09'
$ws.Range("E27").Value = '  +2.80%  '

# Row 28
$ws.Range("E28").Value = '  +0.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.16'
$ws.Range("E29").Value = '  +7.52%  '

# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.58'
$ws.Range("E30").Value = '  +2.04%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.15'
$ws.Range("E31").Value = '  +4.74%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0000100'
$ws.Range("E32").Value = '  +18.85%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.108'
$ws.Range("E33").Value = '  +4.08%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.13'
$ws.Range("E34").Value = '  +2.56%  '

# Row 35
$ws.Range("E35").Value = '  -0.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.991'
$ws.Range("E36").Value = '  +2.90%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.57'
$ws.Range("E37").Value = '  +2.77%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.00'
$ws.Range("E38").Value = '  +7.68%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.64'
$ws.Range("E39").Value = '  +1.53%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("E40").Value = '  +6.43%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.33'
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.114'
$ws.Range("E42").Value = '  -2.22%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.273'
$ws.Range("E43").Value = '  +4.28%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.56'
$ws.Range("E44").Value = '  -0.26%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '134.99'
$ws.Range("E45").Value = '  +2.11%  '

# Row 46
$ws.Range("D46").Value = '2.684.60'
$ws.Range("E46").Value = '  +1.25%  '

# Row 47
$ws.Range("E47").Value = '  +1.85%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '354.26'
$ws.Range("E48").Value = '  +3.78%  '

# Row 50
$ws.Range("E50").Value = '  +2.36%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.49'
$ws.Range("E51").Value = '  +0.73%  '
